$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 19 (pushes existing rows 19-35 down to 20-36)
$ws.Rows("19:19").Insert()

# 2. Copy formatting (styles/number formats/borders/etc.) from row 18 into the
#    freshly inserted row 19 so the new line matches the sheet's established
#    look for a data row.
$ws.Range("A18:Q18").Copy()
$ws.Range("A19:Q19").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3. Match row height used by other data rows.
$ws.Rows("19:19").RowHeight = 25.5

# 4. Re-create the merges for the new row (mirrors the pattern used by every
#    other item row: A:B, C:G, H:K, L:M, N:O).
$ws.Range("A19:B19").Merge()
$ws.Range("C19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("N19:O19").Merge()

# 5. Populate the new row with the new item ("LEVOTAVIN 750 MG 7 F.C.TABS.").
$ws.Range("A19").Value2 = 13
$ws.Range("C19").Value2 = "LEVOTAVIN 750 MG 7 F.C.TABS."
$ws.Range("H19").Value2 = "0:0"
$ws.Range("L19").Value2 = "1"
$ws.Range("N19").Value2 = "108.00"
$ws.Range("P19").Value2 = "108.0000"
$ws.Range("Q19").Value2 = "1:0"

# 6. Update the running total (old row 34 -> now row 35) to include the new
#    item's price.
$ws.Range("P35").Value2 = 997.68

# 7. Update the generated-on timestamp (old row 35 -> now row 36).
$ws.Range("A36").Value2 = "Wednesday, 20 August, 2025 4:06 PM"
